$wb = $excel.ActiveWorkbook

# ---- Sheet1 (Step1_Data) and Sheet2 (Step2_Sj): add new column AJ ("Signal_Value_123") ----
$ws1 = $wb.Worksheets.Item("Step1_Data")
$ws2 = $wb.Worksheets.Item("Step2_Sj")

foreach ($ws in @($ws1, $ws2)) {
    $ws.Cells.Item(1, 36).Value = "Signal_Value_123"
    $ws.Cells.Item(1, 36).Font.Bold = $true
    $ws.Cells.Item(1, 36).HorizontalAlignment = -4108
    $ws.Cells.Item(1, 36).VerticalAlignment = -4160
    $ws.Cells.Item(1, 36).Borders.LineStyle = 1
}

# -- Sheet1 raw signal values (rows 2-6, cols B..AJ) --
# row 2: 'signal segment 1'
$ws1.Cells.Item(2, 2).Value = 0
$ws1.Cells.Item(2, 3).Value = 0
$ws1.Cells.Item(2, 4).Value = 0
$ws1.Cells.Item(2, 5).Value = 0.2292832189136331
$ws1.Cells.Item(2, 6).Value = 0.09407438916884261
$ws1.Cells.Item(2, 7).Value = 0.1110127260239029
$ws1.Cells.Item(2, 8).Value = 0.0004462477610108365
$ws1.Cells.Item(2, 9).Value = 0
$ws1.Cells.Item(2, 10).Value = 0
$ws1.Cells.Item(2, 11).Value = 0
$ws1.Cells.Item(2, 12).Value = 0.01567226188855268
$ws1.Cells.Item(2, 13).Value = 0.04787201471968249
$ws1.Cells.Item(2, 14).Value = 0
$ws1.Cells.Item(2, 15).Value = 0.1452884344364642
$ws1.Cells.Item(2, 16).Value = 0
$ws1.Cells.Item(2, 17).Value = 0
$ws1.Cells.Item(2, 18).Value = 0.008452096457880462
$ws1.Cells.Item(2, 19).Value = 0.01978881127667173
$ws1.Cells.Item(2, 20).Value = 0.05413820201975984
$ws1.Cells.Item(2, 21).Value = 0.01366672262426629
$ws1.Cells.Item(2, 22).Value = 0
$ws1.Cells.Item(2, 23).Value = 0.07006685818470805
$ws1.Cells.Item(2, 24).Value = 0
$ws1.Cells.Item(2, 25).Value = 0.03442990511601708
$ws1.Cells.Item(2, 26).Value = 0.0069178816625737
$ws1.Cells.Item(2, 27).Value = 0.087221621858022
$ws1.Cells.Item(2, 28).Value = 0
$ws1.Cells.Item(2, 29).Value = 0.01915629669121869
$ws1.Cells.Item(2, 30).Value = 0.004558322243269392
$ws1.Cells.Item(2, 31).Value = 0
$ws1.Cells.Item(2, 32).Value = 0.03068981887159912
$ws1.Cells.Item(2, 33).Value = 0
$ws1.Cells.Item(2, 34).Value = 0
$ws1.Cells.Item(2, 35).Value = 0
$ws1.Cells.Item(2, 36).Value = 0.00726417008192482
# row 3: 'signal segment 2'
$ws1.Cells.Item(3, 2).Value = 0
$ws1.Cells.Item(3, 3).Value = 0
$ws1.Cells.Item(3, 4).Value = 0.1129906515350862
$ws1.Cells.Item(3, 5).Value = 0
$ws1.Cells.Item(3, 6).Value = 0.3248919432741574
$ws1.Cells.Item(3, 7).Value = 0
$ws1.Cells.Item(3, 8).Value = 0
$ws1.Cells.Item(3, 9).Value = 0
$ws1.Cells.Item(3, 10).Value = 0
$ws1.Cells.Item(3, 11).Value = 0.002122499761286698
$ws1.Cells.Item(3, 12).Value = 0
$ws1.Cells.Item(3, 13).Value = 0.02670598144566798
$ws1.Cells.Item(3, 14).Value = 0.07223171127067399
$ws1.Cells.Item(3, 15).Value = 0.07478170243417812
$ws1.Cells.Item(3, 16).Value = 0
$ws1.Cells.Item(3, 17).Value = 0.03510610811339969
$ws1.Cells.Item(3, 18).Value = 0
$ws1.Cells.Item(3, 19).Value = 0.06089293136036755
$ws1.Cells.Item(3, 20).Value = 0
$ws1.Cells.Item(3, 21).Value = 0.001174605266816657
$ws1.Cells.Item(3, 22).Value = 0.0130805804415602
$ws1.Cells.Item(3, 23).Value = 0.02149834292715637
$ws1.Cells.Item(3, 24).Value = 0.04323958725855117
$ws1.Cells.Item(3, 25).Value = 0
$ws1.Cells.Item(3, 26).Value = 0.09596772067934939
$ws1.Cells.Item(3, 27).Value = 0.01580263954144162
$ws1.Cells.Item(3, 28).Value = 0.02096797666967359
$ws1.Cells.Item(3, 29).Value = 0
$ws1.Cells.Item(3, 30).Value = 0.008440932476203198
$ws1.Cells.Item(3, 31).Value = 0.01428927318480741
$ws1.Cells.Item(3, 32).Value = 0.03106330025600925
$ws1.Cells.Item(3, 33).Value = 0
$ws1.Cells.Item(3, 34).Value = 0
$ws1.Cells.Item(3, 35).Value = 0.02475151210361332
$ws1.Cells.Item(3, 36).Value = 0
# row 4: 'signal segment 3'
$ws1.Cells.Item(4, 2).Value = 0
$ws1.Cells.Item(4, 3).Value = 0
$ws1.Cells.Item(4, 4).Value = 0.193986046549966
$ws1.Cells.Item(4, 5).Value = 0.1786438391707103
$ws1.Cells.Item(4, 6).Value = 0.09493882529734858
$ws1.Cells.Item(4, 7).Value = 0.0129459521539048
$ws1.Cells.Item(4, 8).Value = 0
$ws1.Cells.Item(4, 9).Value = 0.002346011367516524
$ws1.Cells.Item(4, 10).Value = 0
$ws1.Cells.Item(4, 11).Value = 0.011625756205607
$ws1.Cells.Item(4, 12).Value = 0.06344820762750607
$ws1.Cells.Item(4, 13).Value = 0.001021357383045717
$ws1.Cells.Item(4, 14).Value = 0.1238766761379425
$ws1.Cells.Item(4, 15).Value = 0
$ws1.Cells.Item(4, 16).Value = 0
$ws1.Cells.Item(4, 17).Value = 0.002868291652127786
$ws1.Cells.Item(4, 18).Value = 0.02537433479874659
$ws1.Cells.Item(4, 19).Value = 0.04029084771273393
$ws1.Cells.Item(4, 20).Value = 0.02400465862045287
$ws1.Cells.Item(4, 21).Value = 0
$ws1.Cells.Item(4, 22).Value = 0.0671086491683819
$ws1.Cells.Item(4, 23).Value = 0
$ws1.Cells.Item(4, 24).Value = 0.01999699497264572
$ws1.Cells.Item(4, 25).Value = 0.02746254083177083
$ws1.Cells.Item(4, 26).Value = 0.04364875010612293
$ws1.Cells.Item(4, 27).Value = 0
$ws1.Cells.Item(4, 28).Value = 0.01009751856923567
$ws1.Cells.Item(4, 29).Value = 0.01238486842080773
$ws1.Cells.Item(4, 30).Value = 0
$ws1.Cells.Item(4, 31).Value = 0.03260814882657146
$ws1.Cells.Item(4, 32).Value = 0
$ws1.Cells.Item(4, 33).Value = 0
$ws1.Cells.Item(4, 34).Value = 0.01132172442685512
$ws1.Cells.Item(4, 35).Value = 0
$ws1.Cells.Item(4, 36).Value = 0
# row 5: 'signal segment 4'
$ws1.Cells.Item(5, 2).Value = 0
$ws1.Cells.Item(5, 3).Value = 0
$ws1.Cells.Item(5, 4).Value = 0.1988540320646729
$ws1.Cells.Item(5, 5).Value = 0.06965070396710306
$ws1.Cells.Item(5, 6).Value = 0.1358158410711957
$ws1.Cells.Item(5, 7).Value = 0.006152263830569653
$ws1.Cells.Item(5, 8).Value = 0
$ws1.Cells.Item(5, 9).Value = 0
$ws1.Cells.Item(5, 10).Value = 0
$ws1.Cells.Item(5, 11).Value = 0.01170190227478735
$ws1.Cells.Item(5, 12).Value = 0.04201164203131198
$ws1.Cells.Item(5, 13).Value = 0
$ws1.Cells.Item(5, 14).Value = 0.1473342503636893
$ws1.Cells.Item(5, 15).Value = 0.007383888863676438
$ws1.Cells.Item(5, 16).Value = 0
$ws1.Cells.Item(5, 17).Value = 0.01831648096308858
$ws1.Cells.Item(5, 18).Value = 0.02117276378901272
$ws1.Cells.Item(5, 19).Value = 0.06755243545549391
$ws1.Cells.Item(5, 20).Value = 0.0160492635212033
$ws1.Cells.Item(5, 21).Value = 0
$ws1.Cells.Item(5, 22).Value = 0.05982115054705688
$ws1.Cells.Item(5, 23).Value = 0
$ws1.Cells.Item(5, 24).Value = 0.03269358441674861
$ws1.Cells.Item(5, 25).Value = 0.0004614332732031369
$ws1.Cells.Item(5, 26).Value = 0.08813558481449985
$ws1.Cells.Item(5, 27).Value = 0
$ws1.Cells.Item(5, 28).Value = 0.02454386834303442
$ws1.Cells.Item(5, 29).Value = 0.0003184417548803606
$ws1.Cells.Item(5, 30).Value = 0
$ws1.Cells.Item(5, 31).Value = 0.0355425654284788
$ws1.Cells.Item(5, 32).Value = 0.003225033501748069
$ws1.Cells.Item(5, 33).Value = 0
$ws1.Cells.Item(5, 34).Value = 0
$ws1.Cells.Item(5, 35).Value = 0.01326286972454487
$ws1.Cells.Item(5, 36).Value = 0
# row 6: 'signal segment 5'
$ws1.Cells.Item(6, 2).Value = 0
$ws1.Cells.Item(6, 3).Value = 0
$ws1.Cells.Item(6, 4).Value = 0.186314361635917
$ws1.Cells.Item(6, 5).Value = 0.1783129619679877
$ws1.Cells.Item(6, 6).Value = 0.05581812392977374
$ws1.Cells.Item(6, 7).Value = 0
$ws1.Cells.Item(6, 8).Value = 0
$ws1.Cells.Item(6, 9).Value = 0
$ws1.Cells.Item(6, 10).Value = 0
$ws1.Cells.Item(6, 11).Value = 0
$ws1.Cells.Item(6, 12).Value = 0.04896967811916524
$ws1.Cells.Item(6, 13).Value = 0.006804712120033768
$ws1.Cells.Item(6, 14).Value = 0.1399804087345639
$ws1.Cells.Item(6, 15).Value = 0
$ws1.Cells.Item(6, 16).Value = 0
$ws1.Cells.Item(6, 17).Value = 0
$ws1.Cells.Item(6, 18).Value = 0.03157084956095623
$ws1.Cells.Item(6, 19).Value = 0.02628740681755598
$ws1.Cells.Item(6, 20).Value = 0.01726131133512003
$ws1.Cells.Item(6, 21).Value = 0
$ws1.Cells.Item(6, 22).Value = 0.07536033535755536
$ws1.Cells.Item(6, 23).Value = 0
$ws1.Cells.Item(6, 24).Value = 0.03443072970759244
$ws1.Cells.Item(6, 25).Value = 0.02239987436689955
$ws1.Cells.Item(6, 26).Value = 0.09098122660098146
$ws1.Cells.Item(6, 27).Value = 0
$ws1.Cells.Item(6, 28).Value = 0.02209139441855666
$ws1.Cells.Item(6, 29).Value = 0.01438170729951912
$ws1.Cells.Item(6, 30).Value = 0
$ws1.Cells.Item(6, 31).Value = 0.03816866043873077
$ws1.Cells.Item(6, 32).Value = 0
$ws1.Cells.Item(6, 33).Value = 0
$ws1.Cells.Item(6, 34).Value = 0.005842581653868572
$ws1.Cells.Item(6, 35).Value = 0.005023675935222437
$ws1.Cells.Item(6, 36).Value = 0

# -- Sheet2 cumulative values (rows 2-6, cols B..AJ) = running total of Sheet1 from D onward --
# row 2: 'signal segment 1'
$ws2.Cells.Item(2, 2).Value = 0.0
$ws2.Cells.Item(2, 3).Value = 0.0
$ws2.Cells.Item(2, 4).Value = 0.0
$ws2.Cells.Item(2, 5).Value = 0.2292832189136331
$ws2.Cells.Item(2, 6).Value = 0.3233576080824757
$ws2.Cells.Item(2, 7).Value = 0.4343703341063786
$ws2.Cells.Item(2, 8).Value = 0.4348165818673894
$ws2.Cells.Item(2, 9).Value = 0.4348165818673894
$ws2.Cells.Item(2, 10).Value = 0.4348165818673894
$ws2.Cells.Item(2, 11).Value = 0.4348165818673894
$ws2.Cells.Item(2, 12).Value = 0.4504888437559421
$ws2.Cells.Item(2, 13).Value = 0.4983608584756246
$ws2.Cells.Item(2, 14).Value = 0.4983608584756246
$ws2.Cells.Item(2, 15).Value = 0.6436492929120888
$ws2.Cells.Item(2, 16).Value = 0.6436492929120888
$ws2.Cells.Item(2, 17).Value = 0.6436492929120888
$ws2.Cells.Item(2, 18).Value = 0.6521013893699693
$ws2.Cells.Item(2, 19).Value = 0.6718902006466411
$ws2.Cells.Item(2, 20).Value = 0.7260284026664009
$ws2.Cells.Item(2, 21).Value = 0.7396951252906672
$ws2.Cells.Item(2, 22).Value = 0.7396951252906672
$ws2.Cells.Item(2, 23).Value = 0.8097619834753752
$ws2.Cells.Item(2, 24).Value = 0.8097619834753752
$ws2.Cells.Item(2, 25).Value = 0.8441918885913923
$ws2.Cells.Item(2, 26).Value = 0.851109770253966
$ws2.Cells.Item(2, 27).Value = 0.9383313921119879
$ws2.Cells.Item(2, 28).Value = 0.9383313921119879
$ws2.Cells.Item(2, 29).Value = 0.9574876888032066
$ws2.Cells.Item(2, 30).Value = 0.9620460110464759
$ws2.Cells.Item(2, 31).Value = 0.9620460110464759
$ws2.Cells.Item(2, 32).Value = 0.9927358299180751
$ws2.Cells.Item(2, 33).Value = 0.9927358299180751
$ws2.Cells.Item(2, 34).Value = 0.9927358299180751
$ws2.Cells.Item(2, 35).Value = 0.9927358299180751
$ws2.Cells.Item(2, 36).Value = 0.9999999999999999
# row 3: 'signal segment 2'
$ws2.Cells.Item(3, 2).Value = 0.0
$ws2.Cells.Item(3, 3).Value = 0.0
$ws2.Cells.Item(3, 4).Value = 0.1129906515350862
$ws2.Cells.Item(3, 5).Value = 0.1129906515350862
$ws2.Cells.Item(3, 6).Value = 0.4378825948092436
$ws2.Cells.Item(3, 7).Value = 0.4378825948092436
$ws2.Cells.Item(3, 8).Value = 0.4378825948092436
$ws2.Cells.Item(3, 9).Value = 0.4378825948092436
$ws2.Cells.Item(3, 10).Value = 0.4378825948092436
$ws2.Cells.Item(3, 11).Value = 0.4400050945705303
$ws2.Cells.Item(3, 12).Value = 0.4400050945705303
$ws2.Cells.Item(3, 13).Value = 0.46671107601619827
$ws2.Cells.Item(3, 14).Value = 0.5389427872868723
$ws2.Cells.Item(3, 15).Value = 0.6137244897210504
$ws2.Cells.Item(3, 16).Value = 0.6137244897210504
$ws2.Cells.Item(3, 17).Value = 0.6488305978344501
$ws2.Cells.Item(3, 18).Value = 0.6488305978344501
$ws2.Cells.Item(3, 19).Value = 0.7097235291948176
$ws2.Cells.Item(3, 20).Value = 0.7097235291948176
$ws2.Cells.Item(3, 21).Value = 0.7108981344616343
$ws2.Cells.Item(3, 22).Value = 0.7239787149031945
$ws2.Cells.Item(3, 23).Value = 0.7454770578303509
$ws2.Cells.Item(3, 24).Value = 0.788716645088902
$ws2.Cells.Item(3, 25).Value = 0.788716645088902
$ws2.Cells.Item(3, 26).Value = 0.8846843657682515
$ws2.Cells.Item(3, 27).Value = 0.9004870053096932
$ws2.Cells.Item(3, 28).Value = 0.9214549819793667
$ws2.Cells.Item(3, 29).Value = 0.9214549819793667
$ws2.Cells.Item(3, 30).Value = 0.9298959144555698
$ws2.Cells.Item(3, 31).Value = 0.9441851876403773
$ws2.Cells.Item(3, 32).Value = 0.9752484878963865
$ws2.Cells.Item(3, 33).Value = 0.9752484878963865
$ws2.Cells.Item(3, 34).Value = 0.9752484878963865
$ws2.Cells.Item(3, 35).Value = 0.9999999999999999
$ws2.Cells.Item(3, 36).Value = 0.9999999999999999
# row 4: 'signal segment 3'
$ws2.Cells.Item(4, 2).Value = 0.0
$ws2.Cells.Item(4, 3).Value = 0.0
$ws2.Cells.Item(4, 4).Value = 0.193986046549966
$ws2.Cells.Item(4, 5).Value = 0.37262988572067635
$ws2.Cells.Item(4, 6).Value = 0.46756871101802494
$ws2.Cells.Item(4, 7).Value = 0.48051466317192976
$ws2.Cells.Item(4, 8).Value = 0.48051466317192976
$ws2.Cells.Item(4, 9).Value = 0.48286067453944626
$ws2.Cells.Item(4, 10).Value = 0.48286067453944626
$ws2.Cells.Item(4, 11).Value = 0.4944864307450533
$ws2.Cells.Item(4, 12).Value = 0.5579346383725593
$ws2.Cells.Item(4, 13).Value = 0.5589559957556051
$ws2.Cells.Item(4, 14).Value = 0.6828326718935476
$ws2.Cells.Item(4, 15).Value = 0.6828326718935476
$ws2.Cells.Item(4, 16).Value = 0.6828326718935476
$ws2.Cells.Item(4, 17).Value = 0.6857009635456753
$ws2.Cells.Item(4, 18).Value = 0.711075298344422
$ws2.Cells.Item(4, 19).Value = 0.7513661460571559
$ws2.Cells.Item(4, 20).Value = 0.7753708046776088
$ws2.Cells.Item(4, 21).Value = 0.7753708046776088
$ws2.Cells.Item(4, 22).Value = 0.8424794538459907
$ws2.Cells.Item(4, 23).Value = 0.8424794538459907
$ws2.Cells.Item(4, 24).Value = 0.8624764488186364
$ws2.Cells.Item(4, 25).Value = 0.8899389896504072
$ws2.Cells.Item(4, 26).Value = 0.9335877397565301
$ws2.Cells.Item(4, 27).Value = 0.9335877397565301
$ws2.Cells.Item(4, 28).Value = 0.9436852583257658
$ws2.Cells.Item(4, 29).Value = 0.9560701267465735
$ws2.Cells.Item(4, 30).Value = 0.9560701267465735
$ws2.Cells.Item(4, 31).Value = 0.988678275573145
$ws2.Cells.Item(4, 32).Value = 0.988678275573145
$ws2.Cells.Item(4, 33).Value = 0.988678275573145
$ws2.Cells.Item(4, 34).Value = 1.0000000000000002
$ws2.Cells.Item(4, 35).Value = 1.0000000000000002
$ws2.Cells.Item(4, 36).Value = 1.0000000000000002
# row 5: 'signal segment 4'
$ws2.Cells.Item(5, 2).Value = 0.0
$ws2.Cells.Item(5, 3).Value = 0.0
$ws2.Cells.Item(5, 4).Value = 0.1988540320646729
$ws2.Cells.Item(5, 5).Value = 0.26850473603177594
$ws2.Cells.Item(5, 6).Value = 0.40432057710297165
$ws2.Cells.Item(5, 7).Value = 0.4104728409335413
$ws2.Cells.Item(5, 8).Value = 0.4104728409335413
$ws2.Cells.Item(5, 9).Value = 0.4104728409335413
$ws2.Cells.Item(5, 10).Value = 0.4104728409335413
$ws2.Cells.Item(5, 11).Value = 0.42217474320832865
$ws2.Cells.Item(5, 12).Value = 0.4641863852396406
$ws2.Cells.Item(5, 13).Value = 0.4641863852396406
$ws2.Cells.Item(5, 14).Value = 0.6115206356033299
$ws2.Cells.Item(5, 15).Value = 0.6189045244670064
$ws2.Cells.Item(5, 16).Value = 0.6189045244670064
$ws2.Cells.Item(5, 17).Value = 0.637221005430095
$ws2.Cells.Item(5, 18).Value = 0.6583937692191076
$ws2.Cells.Item(5, 19).Value = 0.7259462046746016
$ws2.Cells.Item(5, 20).Value = 0.7419954681958049
$ws2.Cells.Item(5, 21).Value = 0.7419954681958049
$ws2.Cells.Item(5, 22).Value = 0.8018166187428617
$ws2.Cells.Item(5, 23).Value = 0.8018166187428617
$ws2.Cells.Item(5, 24).Value = 0.8345102031596103
$ws2.Cells.Item(5, 25).Value = 0.8349716364328135
$ws2.Cells.Item(5, 26).Value = 0.9231072212473133
$ws2.Cells.Item(5, 27).Value = 0.9231072212473133
$ws2.Cells.Item(5, 28).Value = 0.9476510895903477
$ws2.Cells.Item(5, 29).Value = 0.9479695313452281
$ws2.Cells.Item(5, 30).Value = 0.9479695313452281
$ws2.Cells.Item(5, 31).Value = 0.9835120967737069
$ws2.Cells.Item(5, 32).Value = 0.986737130275455
$ws2.Cells.Item(5, 33).Value = 0.986737130275455
$ws2.Cells.Item(5, 34).Value = 0.986737130275455
$ws2.Cells.Item(5, 35).Value = 0.9999999999999998
$ws2.Cells.Item(5, 36).Value = 0.9999999999999998
# row 6: 'signal segment 5'
$ws2.Cells.Item(6, 2).Value = 0.0
$ws2.Cells.Item(6, 3).Value = 0.0
$ws2.Cells.Item(6, 4).Value = 0.186314361635917
$ws2.Cells.Item(6, 5).Value = 0.3646273236039047
$ws2.Cells.Item(6, 6).Value = 0.42044544753367846
$ws2.Cells.Item(6, 7).Value = 0.42044544753367846
$ws2.Cells.Item(6, 8).Value = 0.42044544753367846
$ws2.Cells.Item(6, 9).Value = 0.42044544753367846
$ws2.Cells.Item(6, 10).Value = 0.42044544753367846
$ws2.Cells.Item(6, 11).Value = 0.42044544753367846
$ws2.Cells.Item(6, 12).Value = 0.4694151256528437
$ws2.Cells.Item(6, 13).Value = 0.47621983777287746
$ws2.Cells.Item(6, 14).Value = 0.6162002465074414
$ws2.Cells.Item(6, 15).Value = 0.6162002465074414
$ws2.Cells.Item(6, 16).Value = 0.6162002465074414
$ws2.Cells.Item(6, 17).Value = 0.6162002465074414
$ws2.Cells.Item(6, 18).Value = 0.6477710960683977
$ws2.Cells.Item(6, 19).Value = 0.6740585028859536
$ws2.Cells.Item(6, 20).Value = 0.6913198142210736
$ws2.Cells.Item(6, 21).Value = 0.6913198142210736
$ws2.Cells.Item(6, 22).Value = 0.766680149578629
$ws2.Cells.Item(6, 23).Value = 0.766680149578629
$ws2.Cells.Item(6, 24).Value = 0.8011108792862214
$ws2.Cells.Item(6, 25).Value = 0.8235107536531209
$ws2.Cells.Item(6, 26).Value = 0.9144919802541023
$ws2.Cells.Item(6, 27).Value = 0.9144919802541023
$ws2.Cells.Item(6, 28).Value = 0.936583374672659
$ws2.Cells.Item(6, 29).Value = 0.9509650819721781
$ws2.Cells.Item(6, 30).Value = 0.9509650819721781
$ws2.Cells.Item(6, 31).Value = 0.9891337424109089
$ws2.Cells.Item(6, 32).Value = 0.9891337424109089
$ws2.Cells.Item(6, 33).Value = 0.9891337424109089
$ws2.Cells.Item(6, 34).Value = 0.9949763240647774
$ws2.Cells.Item(6, 35).Value = 0.9999999999999999
$ws2.Cells.Item(6, 36).Value = 0.9999999999999999

# ---- Sheets 3-6 (Step3_DataPts_*): threshold-crossing summaries ----
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Cells.Item(2, 2).Value = 0.5
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = 14
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0.6436492929120887
$ws.Cells.Item(2, 7).Value = 12
$ws.Cells.Item(2, 8).Value = 61
$ws.Cells.Item(2, 9).Value = 700
$ws.Cells.Item(2, 10).Value = "11R22.5"
$ws.Cells.Item(2, 11).Value = "710R"
$ws.Cells.Item(2, 12).Value = "100%"
$ws.Cells.Item(2, 13).Value = "Ir"
$ws.Cells.Item(3, 2).Value = 0.5
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 13
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0.5389427872868723
$ws.Cells.Item(3, 7).Value = 12
$ws.Cells.Item(3, 8).Value = 61
$ws.Cells.Item(3, 9).Value = 700
$ws.Cells.Item(3, 10).Value = "11R22.5"
$ws.Cells.Item(3, 11).Value = "710R"
$ws.Cells.Item(3, 12).Value = "100%"
$ws.Cells.Item(3, 13).Value = "Ir"
$ws.Cells.Item(4, 2).Value = 0.5
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 11
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0.5579346383725592
$ws.Cells.Item(4, 7).Value = 10
$ws.Cells.Item(4, 8).Value = 61
$ws.Cells.Item(4, 9).Value = 700
$ws.Cells.Item(4, 10).Value = "11R22.5"
$ws.Cells.Item(4, 11).Value = "710R"
$ws.Cells.Item(4, 12).Value = "100%"
$ws.Cells.Item(4, 13).Value = "Ir"
$ws.Cells.Item(5, 2).Value = 0.5
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 13
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0.6115206356033299
$ws.Cells.Item(5, 7).Value = 12
$ws.Cells.Item(5, 8).Value = 61
$ws.Cells.Item(5, 9).Value = 700
$ws.Cells.Item(5, 10).Value = "11R22.5"
$ws.Cells.Item(5, 11).Value = "710R"
$ws.Cells.Item(5, 12).Value = "100%"
$ws.Cells.Item(5, 13).Value = "Ir"
$ws.Cells.Item(6, 2).Value = 0.5
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 13
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0.6162002465074413
$ws.Cells.Item(6, 7).Value = 12
$ws.Cells.Item(6, 8).Value = 61
$ws.Cells.Item(6, 9).Value = 700
$ws.Cells.Item(6, 10).Value = "11R22.5"
$ws.Cells.Item(6, 11).Value = "710R"
$ws.Cells.Item(6, 12).Value = "100%"
$ws.Cells.Item(6, 13).Value = "Ir"

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Cells.Item(2, 2).Value = 0.7
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = 19
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0.7260284026664008
$ws.Cells.Item(2, 7).Value = 17
$ws.Cells.Item(2, 8).Value = 61
$ws.Cells.Item(2, 9).Value = 700
$ws.Cells.Item(2, 10).Value = "11R22.5"
$ws.Cells.Item(2, 11).Value = "710R"
$ws.Cells.Item(2, 12).Value = "100%"
$ws.Cells.Item(2, 13).Value = "Ir"
$ws.Cells.Item(3, 2).Value = 0.7
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 18
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0.7097235291948176
$ws.Cells.Item(3, 7).Value = 17
$ws.Cells.Item(3, 8).Value = 61
$ws.Cells.Item(3, 9).Value = 700
$ws.Cells.Item(3, 10).Value = "11R22.5"
$ws.Cells.Item(3, 11).Value = "710R"
$ws.Cells.Item(3, 12).Value = "100%"
$ws.Cells.Item(3, 13).Value = "Ir"
$ws.Cells.Item(4, 2).Value = 0.7
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 17
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0.7110752983444218
$ws.Cells.Item(4, 7).Value = 16
$ws.Cells.Item(4, 8).Value = 61
$ws.Cells.Item(4, 9).Value = 700
$ws.Cells.Item(4, 10).Value = "11R22.5"
$ws.Cells.Item(4, 11).Value = "710R"
$ws.Cells.Item(4, 12).Value = "100%"
$ws.Cells.Item(4, 13).Value = "Ir"
$ws.Cells.Item(5, 2).Value = 0.7
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 18
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0.7259462046746016
$ws.Cells.Item(5, 7).Value = 17
$ws.Cells.Item(5, 8).Value = 61
$ws.Cells.Item(5, 9).Value = 700
$ws.Cells.Item(5, 10).Value = "11R22.5"
$ws.Cells.Item(5, 11).Value = "710R"
$ws.Cells.Item(5, 12).Value = "100%"
$ws.Cells.Item(5, 13).Value = "Ir"
$ws.Cells.Item(6, 2).Value = 0.7
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 21
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0.766680149578629
$ws.Cells.Item(6, 7).Value = 20
$ws.Cells.Item(6, 8).Value = 61
$ws.Cells.Item(6, 9).Value = 700
$ws.Cells.Item(6, 10).Value = "11R22.5"
$ws.Cells.Item(6, 11).Value = "710R"
$ws.Cells.Item(6, 12).Value = "100%"
$ws.Cells.Item(6, 13).Value = "Ir"

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Cells.Item(2, 2).Value = 0.8
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = 22
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0.8097619834753751
$ws.Cells.Item(2, 7).Value = 20
$ws.Cells.Item(2, 8).Value = 61
$ws.Cells.Item(2, 9).Value = 700
$ws.Cells.Item(2, 10).Value = "11R22.5"
$ws.Cells.Item(2, 11).Value = "710R"
$ws.Cells.Item(2, 12).Value = "100%"
$ws.Cells.Item(2, 13).Value = "Ir"
$ws.Cells.Item(3, 2).Value = 0.8
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 25
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0.8846843657682515
$ws.Cells.Item(3, 7).Value = 24
$ws.Cells.Item(3, 8).Value = 61
$ws.Cells.Item(3, 9).Value = 700
$ws.Cells.Item(3, 10).Value = "11R22.5"
$ws.Cells.Item(3, 11).Value = "710R"
$ws.Cells.Item(3, 12).Value = "100%"
$ws.Cells.Item(3, 13).Value = "Ir"
$ws.Cells.Item(4, 2).Value = 0.8
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 21
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0.8424794538459905
$ws.Cells.Item(4, 7).Value = 20
$ws.Cells.Item(4, 8).Value = 61
$ws.Cells.Item(4, 9).Value = 700
$ws.Cells.Item(4, 10).Value = "11R22.5"
$ws.Cells.Item(4, 11).Value = "710R"
$ws.Cells.Item(4, 12).Value = "100%"
$ws.Cells.Item(4, 13).Value = "Ir"
$ws.Cells.Item(5, 2).Value = 0.8
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 21
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0.8018166187428617
$ws.Cells.Item(5, 7).Value = 20
$ws.Cells.Item(5, 8).Value = 61
$ws.Cells.Item(5, 9).Value = 700
$ws.Cells.Item(5, 10).Value = "11R22.5"
$ws.Cells.Item(5, 11).Value = "710R"
$ws.Cells.Item(5, 12).Value = "100%"
$ws.Cells.Item(5, 13).Value = "Ir"
$ws.Cells.Item(6, 2).Value = 0.8
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 23
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0.8011108792862214
$ws.Cells.Item(6, 7).Value = 22
$ws.Cells.Item(6, 8).Value = 61
$ws.Cells.Item(6, 9).Value = 700
$ws.Cells.Item(6, 10).Value = "11R22.5"
$ws.Cells.Item(6, 11).Value = "710R"
$ws.Cells.Item(6, 12).Value = "100%"
$ws.Cells.Item(6, 13).Value = "Ir"

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Cells.Item(2, 2).Value = 0.9
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = 26
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0.9383313921119878
$ws.Cells.Item(2, 7).Value = 24
$ws.Cells.Item(2, 8).Value = 61
$ws.Cells.Item(2, 9).Value = 700
$ws.Cells.Item(2, 10).Value = "11R22.5"
$ws.Cells.Item(2, 11).Value = "710R"
$ws.Cells.Item(2, 12).Value = "100%"
$ws.Cells.Item(2, 13).Value = "Ir"
$ws.Cells.Item(3, 2).Value = 0.9
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 26
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0.9004870053096932
$ws.Cells.Item(3, 7).Value = 25
$ws.Cells.Item(3, 8).Value = 61
$ws.Cells.Item(3, 9).Value = 700
$ws.Cells.Item(3, 10).Value = "11R22.5"
$ws.Cells.Item(3, 11).Value = "710R"
$ws.Cells.Item(3, 12).Value = "100%"
$ws.Cells.Item(3, 13).Value = "Ir"
$ws.Cells.Item(4, 2).Value = 0.9
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 25
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0.93358773975653
$ws.Cells.Item(4, 7).Value = 24
$ws.Cells.Item(4, 8).Value = 61
$ws.Cells.Item(4, 9).Value = 700
$ws.Cells.Item(4, 10).Value = "11R22.5"
$ws.Cells.Item(4, 11).Value = "710R"
$ws.Cells.Item(4, 12).Value = "100%"
$ws.Cells.Item(4, 13).Value = "Ir"
$ws.Cells.Item(5, 2).Value = 0.9
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 25
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0.9231072212473133
$ws.Cells.Item(5, 7).Value = 24
$ws.Cells.Item(5, 8).Value = 61
$ws.Cells.Item(5, 9).Value = 700
$ws.Cells.Item(5, 10).Value = "11R22.5"
$ws.Cells.Item(5, 11).Value = "710R"
$ws.Cells.Item(5, 12).Value = "100%"
$ws.Cells.Item(5, 13).Value = "Ir"
$ws.Cells.Item(6, 2).Value = 0.9
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 25
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0.9144919802541023
$ws.Cells.Item(6, 7).Value = 24
$ws.Cells.Item(6, 8).Value = 61
$ws.Cells.Item(6, 9).Value = 700
$ws.Cells.Item(6, 10).Value = "11R22.5"
$ws.Cells.Item(6, 11).Value = "710R"
$ws.Cells.Item(6, 12).Value = "100%"
$ws.Cells.Item(6, 13).Value = "Ir"

